# Weekly update for "Hortaliza, Vega Central Mapocho de Santiago - Espárragos"
# New price observations were recorded for date 44519 (2021-11-19) at the
# "Provincia de Linares" origin, inserted as the new most-recent rows
# (right after the header block of existing rows), pushing all of the
# subsequent historical rows down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 38:40 - this shifts rows 38..104 down to 41..107
# (preserving their formatting/styles, including the date-formatted column D)
# and grows the sheet's used range from R104 to R107, matching the target.
$ws.Range("A38:R40").Insert()

# Row 38 - Espárragos, Banquete, Provincia de Linares
$ws.Cells.Item(38, 1).Value = 9
$ws.Cells.Item(38, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(38, 3).Value = "Metropolitana"
$ws.Cells.Item(38, 4).Value = 44519
$ws.Cells.Item(38, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(38, 5).Value = 13
$ws.Cells.Item(38, 6).Value = 300000000
$ws.Cells.Item(38, 7).Value = "Espárragos"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Banquete"
$ws.Cells.Item(38, 10).Value = 250
$ws.Cells.Item(38, 11).Value = 1500
$ws.Cells.Item(38, 12).Value = 1500
$ws.Cells.Item(38, 13).Value = 1500
$ws.Cells.Item(38, 14).Value = "`$/kilo"
$ws.Cells.Item(38, 15).Value = "Provincia de Linares"
$ws.Cells.Item(38, 16).Value = 1500
$ws.Cells.Item(38, 17).Value = 1
$ws.Cells.Item(38, 18).Value = "Hortaliza"

# Row 39 - Espárragos, Primera, Provincia de Linares
$ws.Cells.Item(39, 1).Value = 9
$ws.Cells.Item(39, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(39, 3).Value = "Metropolitana"
$ws.Cells.Item(39, 4).Value = 44519
$ws.Cells.Item(39, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(39, 5).Value = 13
$ws.Cells.Item(39, 6).Value = 300000000
$ws.Cells.Item(39, 7).Value = "Espárragos"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 340
$ws.Cells.Item(39, 11).Value = 1300
$ws.Cells.Item(39, 12).Value = 1300
$ws.Cells.Item(39, 13).Value = 1300
$ws.Cells.Item(39, 14).Value = "`$/kilo"
$ws.Cells.Item(39, 15).Value = "Provincia de Linares"
$ws.Cells.Item(39, 16).Value = 1300
$ws.Cells.Item(39, 17).Value = 1
$ws.Cells.Item(39, 18).Value = "Hortaliza"

# Row 40 - Espárragos, Segunda, Provincia de Linares
$ws.Cells.Item(40, 1).Value = 9
$ws.Cells.Item(40, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(40, 3).Value = "Metropolitana"
$ws.Cells.Item(40, 4).Value = 44519
$ws.Cells.Item(40, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(40, 5).Value = 13
$ws.Cells.Item(40, 6).Value = 300000000
$ws.Cells.Item(40, 7).Value = "Espárragos"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Segunda"
$ws.Cells.Item(40, 10).Value = 160
$ws.Cells.Item(40, 11).Value = 1100
$ws.Cells.Item(40, 12).Value = 1100
$ws.Cells.Item(40, 13).Value = 1100
$ws.Cells.Item(40, 14).Value = "`$/kilo"
$ws.Cells.Item(40, 15).Value = "Provincia de Linares"
$ws.Cells.Item(40, 16).Value = 1100
$ws.Cells.Item(40, 17).Value = 1
$ws.Cells.Item(40, 18).Value = "Hortaliza"
